$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("E2").Value = 23.80000000000028
$ws.Range("H2").Value = [double]"3.795634272222758e-16"
$ws.Range("I2").Value = 0.4801349431710304
$ws.Range("K2").Value = 42.77799328627249
$ws.Range("L2").Value = "[36.238362747515815, 49.317623825029166]"
$ws.Range("O2").Value = 1.62897396852804
$ws.Range("P2").Value = "[1.46544762419704, 1.7925003128590404]"
$ws.Range("S2").Value = 59.03593094230053
$ws.Range("T2").Value = "[54.83459760107954, 63.237264283521526]"
$ws.Range("W2").Value = 17.62962962962984
$ws.Range("X2").Value = 17.01021021021041
$ws.Range("Y2").Value = 18.24904904904927

# Row 3 changes
$ws.Range("E3").Value = 23.16000000000018
$ws.Range("G3").Value = [double]"2.642330798607873e-14"
$ws.Range("H3").Value = [double]"8.003365979715415e-14"
$ws.Range("K3").Value = 40.61458504958974
$ws.Range("L3").Value = "[28.177355729575996, 53.05181436960349]"
$ws.Range("M3").Value = [double]"9.546461399168038e-10"
$ws.Range("N3").Value = [double]"9.546461399168038e-10"
$ws.Range("O3").Value = -0.6666843268879239
$ws.Range("P3").Value = "[-0.9685791164220783, -0.3647895373537695]"
$ws.Range("Q3").Value = [double]"2.168433365201672e-05"
$ws.Range("R3").Value = [double]"2.168433365201672e-05"
$ws.Range("S3").Value = 60.8361389464116
$ws.Range("T3").Value = "[54.345414429862146, 67.32686346296106]"
$ws.Range("W3").Value = 2.457417417417439
$ws.Range("X3").Value = 1.344624624624637
$ws.Range("Y3").Value = 3.570210210210241
